$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) / Volume(1h) (column E) values for this refresh of the symbol list.
# Values are stored as text (matching the source sheet), so force a text number format
# before writing them -- otherwise Excel would auto-convert numeric-looking strings like
# "261.14" or "-0.19%" into Number/Percentage cells and normalise their display text.
$updates = @{
    "D2" = "261.14"
    "E2" = "-0.19%"
    "D3" = "27.10"
    "E3" = "-0.65%"
    "D4" = "4.708"
    "E4" = "-0.23%"
    "D5" = "0.06221"
    "E5" = "2.59%"
    "D6" = "6.744"
    "E6" = "0.91%"
    "D7" = "0.8527"
    "E7" = "-1.26%"
    "D8" = "0.9084"
    "E8" = "-1.54%"
    "D9" = "0.1400"
    "E9" = "-0.63%"
    "D10" = "0.04860"
    "E10" = "-4.69%"
    "D11" = "0.07084"
    "E11" = "-1.06%"
    "D12" = "0.03125"
    "E12" = "2.48%"
    "D13" = "0.09058"
    "E13" = "-0.47%"
    "D14" = "0.001530"
    "E14" = "-0.26%"
    "D15" = "0.0006179"
    "E15" = "1.38%"
    "D16" = "0.006121"
    "E16" = "-1.21%"
    "D17" = "3.451"
    "E17" = "0.13%"
    "D18" = "3.174"
    "E18" = "0.13%"
    "E19" = "-0.45%"
    "D22" = "4.090"
    "E22" = "-0.20%"
    "D23" = "0.04245"
    "E23" = "-0.35%"
    "D24" = "0.001213"
    "E24" = "-0.35%"
    "D25" = "0.004087"
    "E25" = "4.48%"
    "E26" = "0.01%"
    "D40" = "0.03905"
    "E40" = "0.56%"
    "E41" = "-0.06%"
    "D42" = "0.004131"
    "E42" = "-0.04%"
    "E43" = "0.99%"
    "D44" = "0.01389"
    "E44" = "-6.85%"
    "D45" = "0.00005121"
    "E45" = "-4.78%"
    "D48" = "0.06691"
    "E48" = "-49.35%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
